$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = 3
$ws.Range("H6").Value = 3
$ws.Range("I6").Value = 4
$ws.Range("J6").Value = 74
$ws.Range("K6").Value = 3

$ws.Range("G7").Value = 6
$ws.Range("H7").Value = 6
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 270
$ws.Range("K7").Value = 4

$ws.Range("G8").Value = 6
$ws.Range("H8").Value = 6
$ws.Range("I8").Value = 2
$ws.Range("J8").Value = 365
$ws.Range("K8").Value = 2

$ws.Range("G9").Value = 6
$ws.Range("H9").Value = 6
$ws.Range("I9").Value = 2
$ws.Range("J9").Value = 234
$ws.Range("K9").Value = 2

$ws.Range("G12").Value = 5
$ws.Range("H12").Value = 5
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 1674
$ws.Range("K12").Value = 5

$ws.Range("G13").Value = 6
$ws.Range("H13").Value = 6
$ws.Range("I13").Value = 0
$ws.Range("J13").Value = 1066
$ws.Range("K13").Value = 4

$ws.Range("G14").Value = 7
$ws.Range("H14").Value = 7
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 691
$ws.Range("K14").Value = 3

$ws.Range("G19").Value = 6
$ws.Range("H19").Value = 6
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 58
$ws.Range("K19").Value = 4

$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 8
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 717
$ws.Range("K20").Value = 2

$ws.Range("G21").Value = 7
$ws.Range("H21").Value = 7
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 181
$ws.Range("K21").Value = 3

$ws.Range("G22").Value = 6
$ws.Range("H22").Value = 6
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 303
$ws.Range("K22").Value = 4

$ws.Range("G24").Value = 6
$ws.Range("H24").Value = 6
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 285
$ws.Range("K24").Value = 4

$ws.Range("G28").Value = 4
$ws.Range("H28").Value = 4
$ws.Range("I28").Value = 0
$ws.Range("J28").Value = 509
$ws.Range("K28").Value = 6

$ws.Range("G29").Value = 4
$ws.Range("H29").Value = 4
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 619
$ws.Range("K29").Value = 6

$ws.Range("G34").Value = 1
$ws.Range("H34").Value = 1
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 50
$ws.Range("K34").Value = 9

$ws.Range("G35").Value = 4
$ws.Range("H35").Value = 4
$ws.Range("I35").Value = 0
$ws.Range("J35").Value = 91
$ws.Range("K35").Value = 6

$ws.Range("G36").Value = 3
$ws.Range("H36").Value = 3
$ws.Range("I36").Value = 0
$ws.Range("J36").Value = 12
$ws.Range("K36").Value = 7

$ws.Range("G42").Value = 9
$ws.Range("H42").Value = 9
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 474
$ws.Range("K42").Value = 1

$ws.Range("G43").Value = 6
$ws.Range("H43").Value = 6
$ws.Range("I43").Value = 0
$ws.Range("J43").Value = 362
$ws.Range("K43").Value = 4

$ws.Range("G44").Value = 7
$ws.Range("H44").Value = 7
$ws.Range("I44").Value = 0
$ws.Range("J44").Value = 273
$ws.Range("K44").Value = 3

$ws.Range("G50").Value = 6
$ws.Range("H50").Value = 6
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 247
$ws.Range("K50").Value = 4

$ws.Range("G51").Value = 5
$ws.Range("H51").Value = 5
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 495
$ws.Range("K51").Value = 5

$ws.Range("G52").Value = 4
$ws.Range("H52").Value = 4
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 319
$ws.Range("K52").Value = 6

$ws.Range("G58").Value = 4
$ws.Range("H58").Value = 4
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 53
$ws.Range("K58").Value = 6

$ws.Range("G59").Value = 4
$ws.Range("H59").Value = 4
$ws.Range("I59").Value = 0
$ws.Range("J59").Value = 80
$ws.Range("K59").Value = 6

$ws.Range("G60").Value = 3
$ws.Range("H60").Value = 3
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 96
$ws.Range("K60").Value = 7

$ws.Range("G63").Value = 0
$ws.Range("H63").Value = 0
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 221
$ws.Range("K63").Value = 10

$ws.Range("G66").Value = 4
$ws.Range("H66").Value = 4
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 1082
$ws.Range("K66").Value = 6

$ws.Range("G67").Value = 3
$ws.Range("H67").Value = 3
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 53
$ws.Range("K67").Value = 7

$ws.Range("G71").Value = 3
$ws.Range("H71").Value = 3
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 183
$ws.Range("K71").Value = 7

$ws.Range("G72").Value = 7
$ws.Range("H72").Value = 7
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 123
$ws.Range("K72").Value = 3

$ws.Range("G73").Value = 4
$ws.Range("H73").Value = 4
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 73
$ws.Range("K73").Value = 6

